# Fix header labels on the existing sheets and add a new "PO Forecast" sheet
# with the Prophet-style forecast output (ds, PO_Forecast, yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# Rename the "Requested quantity" headers to the new metric names.
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new forecast sheet right after "Monthly Trend".
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Header row.
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the bold/bordered/centered header style used on the other sheets.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

# Forecast rows.
$forecastRows = @(
    @(44983.99999999999, 10, 9.99996699935366,  9.999966999715602),
    @(45186.99999999999, 4,  3.999966965251898, 3.999966965616977),
    @(45193.99999999999, 4,  3.793070412325965, 3.793070412713796),
    @(45200.99999999999, 4,  3.586173859413754, 3.586173859817335),
    @(45207.99999999999, 3,  3.379277306588473, 3.379277307035155),
    @(45214.99999999999, 3,  3.172380753630874, 3.172380754143121),
    @(45221.99999999999, 3,  2.965484200651984, 2.965484201218656),
    @(45228.99999999999, 3,  2.758587647719636, 2.758587648363904),
    @(45235.99999999999, 3,  2.551691094102859, 2.551691095572311),
    @(45242.99999999999, 2,  2.344794539641336, 2.34479454381522)
)

$r = 2
foreach ($row in $forecastRows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Match the date-time number format used for the "ds" / "Order Week" column.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A11").PasteSpecial(-4122)

$wsForecast.Range("A1").Select()

Write-Output "Added PO Forecast sheet and renamed headers."
